# Update automatico via Actualizar 02-04-2021 22-40-49
#
# This mirrors what the "Disponibilidad" checker script does on every run:
#   - refresh the timestamp of the most-recent existing check block (rows
#     170-183) to the instant the run actually finished, and
#   - append a brand-new block of 14 rows (one per monitored service) with
#     the new run's timestamp, each with a hyperlink in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Correct the timestamp written for the previous run (rows 170-183).
# ---------------------------------------------------------------------
$refreshedTimestamp = 44231.92382482639
for ($r = 170; $r -le 183; $r++) {
    $ws.Range("D$r").Value2 = $refreshedTimestamp
}

# ---------------------------------------------------------------------
# 2) Append the new run's 14 rows (184-197).
# ---------------------------------------------------------------------
$newTimestamp = 44231.94491828536

$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$displayUrls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$linkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$linkSubAddresses = @($null,$null,$null,$null,$null,$null,$null,$null,"/",$null,$null,$null,$null,$null)

$startRow = 184
for ($i = 0; $i -lt 14; $i++) {
    $row = $startRow + $i

    $ws.Range("A$row").Value2 = $names[$i]
    $ws.Range("B$row").Value2 = $displayUrls[$i]
    $ws.Range("C$row").Value2 = "Disponible"
    $ws.Range("D$row").Value2 = $newTimestamp

    # Match the look of the rest of the sheet: hyperlink-blue/underline style
    # on the URL cell, date/time number format on the timestamp cell.
    $ws.Range("B$row").Style = "Hyperlink"
    $ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $sub = $linkSubAddresses[$i]
    if ($sub) {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $linkAddresses[$i], $sub)
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $linkAddresses[$i])
    }
}
